# Weekly update: insert two new price records (week of 2023-08-09, serial 45147)
# for "Valle de Camiña" ahead of the existing rows, shifting the rest of the
# table down by two rows (old row 586 -> 588, ... old row 616 -> 618).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 586, pushing everything
# from 586..616 down to 588..618.
$ws.Rows.Item(586).Insert()
$ws.Rows.Item(586).Insert()

# New row 586: Primera, Valle de Camiña
$ws.Cells.Item(586, 1).Value = 1
$ws.Cells.Item(586, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(586, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(586, 4).Value = 45147
$ws.Cells.Item(586, 5).Value = 15
$ws.Cells.Item(586, 6).Value = 100114013
$ws.Cells.Item(586, 7).Value = "Zanahoria"
$ws.Cells.Item(586, 8).Value = "Sin especificar"
$ws.Cells.Item(586, 9).Value = "Primera"
$ws.Cells.Item(586, 10).Value = 15
$ws.Cells.Item(586, 11).Value = 12000
$ws.Cells.Item(586, 12).Value = 13000
$ws.Cells.Item(586, 13).Value = 12667
$ws.Cells.Item(586, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(586, 15).Value = "Valle de Camiña"
$ws.Cells.Item(586, 16).Value = 507
$ws.Cells.Item(586, 17).Value = 25
$ws.Cells.Item(586, 18).Value = "Hortaliza"

# New row 587: Segunda, Valle de Camiña
$ws.Cells.Item(587, 1).Value = 1
$ws.Cells.Item(587, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(587, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(587, 4).Value = 45147
$ws.Cells.Item(587, 5).Value = 15
$ws.Cells.Item(587, 6).Value = 100114013
$ws.Cells.Item(587, 7).Value = "Zanahoria"
$ws.Cells.Item(587, 8).Value = "Sin especificar"
$ws.Cells.Item(587, 9).Value = "Segunda"
$ws.Cells.Item(587, 10).Value = 15
$ws.Cells.Item(587, 11).Value = 9000
$ws.Cells.Item(587, 12).Value = 10000
$ws.Cells.Item(587, 13).Value = 9667
$ws.Cells.Item(587, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(587, 15).Value = "Valle de Camiña"
$ws.Cells.Item(587, 16).Value = 387
$ws.Cells.Item(587, 17).Value = 25
$ws.Cells.Item(587, 18).Value = "Hortaliza"
